$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column G: "Deno"
$ws.Range("G3").Value = "Deno"

# Junit row (row 9): bank, Customer, atm -> "y"
$ws.Range("B9").Value = "y"
$ws.Range("C9").Value = "y"
$ws.Range("D9").Value = "y"

# " commits Git" row (row 10): atm, Account, Transaction -> "y"
$ws.Range("D10").Value = "y"
$ws.Range("E10").Value = "y"
$ws.Range("F10").Value = "y"

# Update selection to match the final state (F10)
$ws.Range("F10").Select()
